# Scheduled-runner update: refresh Kraken market-board snapshot values
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2710.3333
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2924.125
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 8772.375
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -9108.375
$ws.Range("H29").Value = 2801
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 3668
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 11004
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -11566
$ws.Range("H40").Value = 9999.75
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9999.75
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9999.75
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -10349.75
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H97").Value = 7165.8335
$ws.Range("J97").Value = 7165.8335
$ws.Range("L97").Value = 21497.5005
$ws.Range("N97").Value = -22489.5005
$ws.Range("H127").Value = 744.3333
$ws.Range("I127").Value = 744.3333
$ws.Range("K127").Value = 2232.9999
$ws.Range("M127").Value = 2727.0001
$ws.Range("H137").Value = 3382.3
$ws.Range("I137").Value = 2611
$ws.Range("J137").Value = 3828.842
$ws.Range("K137").Value = 7833
$ws.Range("L137").Value = 11486.526
$ws.Range("M137").Value = -5283
$ws.Range("N137").Value = -16586.526
$ws.Range("H138").Value = 6264.143
$ws.Range("J138").Value = 6321.75
$ws.Range("L138").Value = 18965.25
$ws.Range("N138").Value = -29245.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 275
$ws.Range("I4").Value = 275
$ws.Range("K4").Value = 275
$ws.Range("M4").Value = -159
$ws.Range("H32").Value = 11995.5625
$ws.Range("I32").Value = 8359.817999999999
$ws.Range("K32").Value = 8359.817999999999
$ws.Range("M32").Value = -8072.817999999999
$ws.Range("H44").Value = 26000.166
$ws.Range("I44").Value = 8000.5
$ws.Range("K44").Value = 8000.5
$ws.Range("M44").Value = -7512.5
$ws.Range("H61").Value = 2833
$ws.Range("I61").Value = 2833
$ws.Range("K61").Value = 2833
$ws.Range("M61").Value = -2621
$ws.Range("H110").Value = 1324.75
$ws.Range("I110").Value = 1324.75
$ws.Range("K110").Value = 1324.75
$ws.Range("M110").Value = 720.25
$ws.Range("H136").Value = 2833
$ws.Range("I136").Value = 2833
$ws.Range("K136").Value = 8499
$ws.Range("M136").Value = -5949

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3619.4707
$ws.Range("I86").Value = 1783.2
$ws.Range("K86").Value = 1783.2
$ws.Range("M86").Value = -660.2
$ws.Range("H89").Value = 3619.4707
$ws.Range("I89").Value = 1783.2
$ws.Range("K89").Value = 8916
$ws.Range("M89").Value = -3300
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 112500750
$ws.Range("J4").Value = 133334340
$ws.Range("L4").Value = 133334340
$ws.Range("N4").Value = -133334564
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H58").Value = 3496.75
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 3500
$ws.Range("N58").Value = -3906
$ws.Range("H132").Value = 2136.6875
$ws.Range("I132").Value = 2118.9167
$ws.Range("K132").Value = 6356.750100000001
$ws.Range("M132").Value = -3826.750100000001
$ws.Range("H134").Value = 1997.5
$ws.Range("I134").Value = 1997
$ws.Range("J134").Value = 1998
$ws.Range("K134").Value = 5991
$ws.Range("L134").Value = 5994
$ws.Range("M134").Value = -3456
$ws.Range("N134").Value = -11064
$ws.Range("H136").Value = 3496.75
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -15600

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H131").Value = 2696.4375
$ws.Range("J131").Value = 3933.3333
$ws.Range("L131").Value = 11799.9999
$ws.Range("N131").Value = -21879.9999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3500
$ws.Range("J80").Value = 3500
$ws.Range("L80").Value = 3500
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 3500
$ws.Range("J83").Value = 3500
$ws.Range("L83").Value = 17500
$ws.Range("N83").Value = -27484
$ws.Range("H97").Value = 3663
$ws.Range("I97").Value = 1994.5
$ws.Range("K97").Value = 1994.5
$ws.Range("M97").Value = -1498.5
$ws.Range("H102").Value = 9166.666999999999
$ws.Range("I102").Value = 9166.666999999999
$ws.Range("K102").Value = 9166.666999999999
$ws.Range("M102").Value = -7544.666999999999
$ws.Range("H134").Value = 39997.5
$ws.Range("J134").Value = 39997.5
$ws.Range("L134").Value = 119992.5
$ws.Range("N134").Value = -125062.5
$ws.Range("H136").Value = 80000
$ws.Range("J136").Value = 80000
$ws.Range("L136").Value = 240000
$ws.Range("N136").Value = -245100
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 15000000
$ws.Range("J2").Value = 15000000
$ws.Range("L2").Value = 15000000
$ws.Range("N2").Value = -15000224
$ws.Range("H21").Value = 20670.334
$ws.Range("I21").Value = 20670.334
$ws.Range("K21").Value = 20670.334
$ws.Range("M21").Value = -20496.334
$ws.Range("H46").Value = 899.8333
$ws.Range("J46").Value = 999.6667
$ws.Range("L46").Value = 999.6667
$ws.Range("N46").Value = -1375.6667
$ws.Range("H55").Value = 6333
$ws.Range("J55").Value = 7500.5
$ws.Range("L55").Value = 7500.5
$ws.Range("N55").Value = -7846.5
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H82").Value = 2328.1428
$ws.Range("I82").Value = 2328.1428
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2328.1428
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1967.1428
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2328.1428
$ws.Range("I85").Value = 2328.1428
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2328.1428
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1080.1428
$ws.Range("N85").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5003426
$ws.Range("I14").Value = 5003426
$ws.Range("K14").Value = 5003426
$ws.Range("M14").Value = -5003258
$ws.Range("H30").Value = 15502
$ws.Range("I30").Value = 20336
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 20336
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -20229
$ws.Range("N30").Value = -1214
$ws.Range("H54").Value = 23552.666
$ws.Range("J54").Value = 34916.5
$ws.Range("L54").Value = 34916.5
$ws.Range("N54").Value = -35956.5
$ws.Range("H122").Value = 502124.75
$ws.Range("I122").Value = 1000500
$ws.Range("K122").Value = 3001500
$ws.Range("M122").Value = -2999050
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H136").Value = 2245.5881
$ws.Range("I136").Value = 2055.4285
$ws.Range("K136").Value = 6166.2855
$ws.Range("M136").Value = -3616.2855
